$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename shared string "specialist" -> "specialistInitial" ---
# This is the text in C7 ("specialist" MBS item), being clarified now that
# a repeat-visit row already exists ("specialistRepeat") and we are about
# to add more rows, so we disambiguate it as the "initial" specialist visit.
$ws.Range("C7").Value() = "specialistInitial"

# --- Add a new cost row for the State surveillance diseases (STEC / L02B) ---
# Copy formatting from an existing "DRG code" style row (row 46) down onto
# the new row 52, then overwrite the values with the new data.
$ws.Range("A46:D46").Copy()
$ws.Range("A52").PasteSpecial(-4122)

$ws.Range("A52").Value() = 6065
$ws.Range("B52").Value() = "L02B"
$ws.Range("C52").Value() = "L02B"
$ws.Range("D52").Value() = "DRG code"

# --- Update the view so the new row is visible / selected ---
$win = $excel.ActiveWindow
$win.ScrollRow = 32
$win.ScrollColumn = 1
$ws.Range("A53").Select()
